$d = $word.ActiveDocument

$xml = "<?xml version=`"1.0`" standalone=`"yes`"?><?mso-application progid=`"Word.Document`"?><pkg:package xmlns:pkg=`"http://schemas.microsoft.com/office/2006/xmlPackage`"><pkg:part pkg:name=`"/word/document.xml`" pkg:contentType=`"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml`"><pkg:xmlData><w:document xmlns:w=`"http://schemas.openxmlformats.org/wordprocessingml/2006/main`"><w:body><w:p><w:r><w:tab/></w:r><w:r><w:t xml:space=`"preserve`">Nowadays, the gradual disappearance of rainforests worldwide has raised the concern of the public, and the last of the rainforests might have permanently ceased to exist within a few decades. </w:t></w:r><w:r><w:t>The extinction under way prompted people to figure out the cause behind the issue and start to consider possible solutions to prevent any unknown catastrophic losses.</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:rFonts w:hint=`"eastAsia`"/></w:rPr></w:pPr><w:r><w:tab/></w:r><w:r><w:t xml:space=`"preserve`">First, there are multiple sources that could lead to the extinction of rainforests, but the most crucial one must be carbon dioxide emissions. </w:t></w:r><w:r><w:t>Moreover, to explore the roots of carbon dioxide emissions, we should discuss why tons of carbon dioxide are being emitted every day, which can be separated into two parts to probe: vehicles and factories.</w:t></w:r><w:r><w:t xml:space=`"preserve`"> </w:t></w:r><w:r><w:t>Automobiles had already become an important part of our daily lives, but using them not only produces tons of carbon dioxide day by day but also causes traffic jams, which can make someone else feel annoying.</w:t></w:r><w:r><w:rPr><w:rFonts w:hint=`"eastAsia`"/></w:rPr><w:t xml:space=`"preserve`"> </w:t></w:r><w:r><w:t>To deal with it, the authorities concerned should make the laws change and give some substitute plans, like buses, MRT, trains, etc., at first.</w:t></w:r><w:r><w:rPr><w:rFonts w:hint=`"eastAsia`"/></w:rPr><w:t xml:space=`"preserve`"> </w:t></w:r><w:r><w:t xml:space=`"preserve`">Subsequently, people should try to accept public transportation and gradually use it in their lives to reduce carbon dioxide emissions. </w:t></w:r></w:p><w:p><w:pPr><w:ind w:firstLine=`"720`"/></w:pPr><w:r><w:t>Secondly, although people use public transportation beforehand, factories might also become murderers of rainforests because they emit an excessive amount of carbon dioxide.</w:t></w:r><w:r><w:t xml:space=`"preserve`"> </w:t></w:r><w:r><w:t>Furthermore, there are two main reasons why factories become murderers: emitting gases without filters and not following eco-friendly procedures.</w:t></w:r><w:r><w:t xml:space=`"preserve`"> </w:t></w:r><w:r><w:t>Nevertheless, although both of them come from a lack of awareness of factories, the legislators should be responsible for them as well due to the loosening of regulations.</w:t></w:r><w:r><w:t xml:space=`"preserve`"> </w:t></w:r><w:r><w:t>Therefore, the urgent priority is to make a far more restrictive regulation to force those violating factories to install the filters, improve the procedures, etc.</w:t></w:r><w:r><w:rPr><w:rFonts w:hint=`"eastAsia`"/></w:rPr><w:t xml:space=`"preserve`"> </w:t></w:r><w:r><w:t xml:space=`"preserve`">Despite </w:t></w:r><w:r><w:t>legislators’</w:t></w:r><w:r><w:t xml:space=`"preserve`"> responsibilities, people can make efforts to prevent it by refusing to use the products produced by those violating factories.</w:t></w:r></w:p><w:p><w:pPr><w:ind w:firstLine=`"720`"/><w:rPr><w:rFonts w:hint=`"eastAsia`"/></w:rPr></w:pPr><w:r><w:t>In conclusion, to prevent those precious species from extinction and create a better world, not only the authorities concerned but also people should make efforts to reduce carbon dioxide emissions significantly from vehicles and factories.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>"

$lastPara = $d.Paragraphs.Last
$insertionPoint = $lastPara.Range
$insertionPoint.Collapse(1)
$insertionPoint.InsertXML($xml)

# The original trailing empty paragraph is still present after the newly
# inserted content; merge it away by deleting the paragraph mark that sits
# between our last inserted paragraph and that now-redundant empty one.
$newLast = $d.Paragraphs.Last
$newPrev = $newLast.Previous
$mergeRange = $d.Range($newPrev.Range.End - 1, $newLast.Range.Start)
$mergeRange.Delete()
